$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing row 15 task description (more detailed commit text)
$ws.Range("B15").Value = "Pom cleaning , add completed featuers to README.ME file, some refactorization"

# Add a new row describing the latest task
$ws.Range("A16").Value = "#100014"
$ws.Range("B16").Value = "Add keycloak logging based on docker container"

# Move the selection to reflect where the user ended up after editing
$ws.Range("B20").Select()
